# AFDP-559 - update Drools "Save Case File" rule text/labels and refresh the
# sheet's saved view (top-left cell / selection) and a couple of row heights.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the rule set / rule table to reflect the CaseFile (not Complaint) domain.
$ws.Range("D2").Value = "Save Case File Rules"
$ws.Range("C12").Value = "RuleTable Save Case File Rules"

# Tighten a couple of row heights that Excel recalculated on save.
$ws.Rows.Item(2).RowHeight = 13.3
$ws.Rows.Item(12).RowHeight = 13.3

# Refresh the saved scroll position / active selection for the sheet.
$aw = $excel.ActiveWindow
$aw.ScrollRow = 1
$aw.ScrollColumn = 1
$ws.Range("C13").Select()
